{"js": "// Update vehicle listing document:\n//  - Ferrari -F8 Triturbo: new Potenza/Cilindrata + new Prezzo\n//  - BMW -S1000RR: new Cilindrata (Potenza/Prezzo unchanged)\n//  - Kawasaki -Ninja -> Yamaha -R1 (new name, Potenza/Cilindrata, Prezzo)\n//  - Porsche -911 S -> Ford -Mustang GT (new name, Potenza/Cilindrata, Prezzo)\n//  - Tesla -Model 3 paragraph removed entirely\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, \"Replace\");\n    await context.sync();\n  }\n}\n\n// 1. Ferrari -F8 Triturbo\nawait replaceOnce(\n  \"   -Potenza: 598,99kw   -Cilindrata: 3500 cm2\",\n  \"   -Potenza: 459kw   -Cilindrata: 4000 cm2\"\n);\nawait replaceOnce(\"  -Prezzo: 380000$\", \"  -Prezzo: 388999$\");\n\n// 2. BMW -S1000RR\nawait replaceOnce(\n  \"   -Potenza: 150kw   -Cilindrata: 1100 cm2\",\n  \"   -Potenza: 150kw   -Cilindrata: 1000 cm2\"\n);\n\n// 3. Kawasaki -Ninja -> Yamaha -R1\nawait replaceOnce(\"Kawasaki -Ninja\", \"Yamaha -R1\");\nawait replaceOnce(\n  \"   -Potenza: 270kw   -Cilindrata: 1200 cm2\",\n  \"   -Potenza: 180kw   -Cilindrata: 1000 cm2\"\n);\nawait replaceOnce(\"  -Prezzo: 40000$\", \"  -Prezzo: 32000$\");\n\n// 4. Porsche -911 S -> Ford -Mustang GT\nawait replaceOnce(\"Porsche -911 S\", \"Ford -Mustang GT\");\nawait replaceOnce(\n  \"   -Potenza: 350kw   -Cilindrata: 3000 cm2\",\n  \"   -Potenza: 400kw   -Cilindrata: 3500 cm2\"\n);\nawait replaceOnce(\"  -Prezzo: 190000$\", \"  -Prezzo: 55000$\");\n\n// 5. Remove the \"Tesla -Model 3\" paragraph entirely (whole block, incl. its\n//    heading run and the MOTORIZZAZIONE/Potenza/Prezzo run).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Tesla -Model 3\") === 0) {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Update vehicle listing document:\n#  - Ferrari -F8 Triturbo: new Potenza/Cilindrata + new Prezzo\n#  - BMW -S1000RR: new Cilindrata (Potenza/Prezzo unchanged)\n#  - Kawasaki -Ninja -> Yamaha -R1 (new name, Potenza/Cilindrata, Prezzo)\n#  - Porsche -911 S -> Ford -Mustang GT (new name, Potenza/Cilindrata, Prezzo)\n#  - Tesla -Model 3 paragraph removed entirely\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceOne = 1 (used for single, precise replacements)\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n}\n\n# 1. Ferrari -F8 Triturbo\nReplace-Text \"   -Potenza: 598,99kw   -Cilindrata: 3500 cm2\" \"   -Potenza: 459kw   -Cilindrata: 4000 cm2\"\nReplace-Text \"  -Prezzo: 380000$\" \"  -Prezzo: 388999$\"\n\n# 2. BMW -S1000RR\nReplace-Text \"   -Potenza: 150kw   -Cilindrata: 1100 cm2\" \"   -Potenza: 150kw   -Cilindrata: 1000 cm2\"\n\n# 3. Kawasaki -Ninja -> Yamaha -R1\nReplace-Text \"Kawasaki -Ninja\" \"Yamaha -R1\"\nReplace-Text \"   -Potenza: 270kw   -Cilindrata: 1200 cm2\" \"   -Potenza: 180kw   -Cilindrata: 1000 cm2\"\nReplace-Text \"  -Prezzo: 40000$\" \"  -Prezzo: 32000$\"\n\n# 4. Porsche -911 S -> Ford -Mustang GT\nReplace-Text \"Porsche -911 S\" \"Ford -Mustang GT\"\nReplace-Text \"   -Potenza: 350kw   -Cilindrata: 3000 cm2\" \"   -Potenza: 400kw   -Cilindrata: 3500 cm2\"\nReplace-Text \"  -Prezzo: 190000$\" \"  -Prezzo: 55000$\"\n\n# 5. Remove the \"Tesla -Model 3\" paragraph entirely (whole block, incl. its\n#    heading run and the MOTORIZZAZIONE/Potenza/Prezzo run).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Tesla -Model 3*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
